$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 261.84616
$ws.Range("I33").Value = 286.94116
$ws.Range("J33").Value = 91.2
$ws.Range("K33").Value = 286.94116
$ws.Range("L33").Value = 91.2
$ws.Range("M33").Value = -57.94116000000002
$ws.Range("N33").Value = -549.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1714.6154
$ws.Range("I40").Value = 1496.6666
$ws.Range("K40").Value = 1496.6666
$ws.Range("M40").Value = -1321.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3252.6667
$ws.Range("I64").Value = 2916.6667
$ws.Range("J64").Value = 3476.6667
$ws.Range("K64").Value = 2916.6667
$ws.Range("L64").Value = 3476.6667
$ws.Range("M64").Value = -2668.6667
$ws.Range("N64").Value = -3972.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3252.6667
$ws.Range("I67").Value = 2916.6667
$ws.Range("J67").Value = 3476.6667
$ws.Range("K67").Value = 2916.6667
$ws.Range("L67").Value = 3476.6667
$ws.Range("M67").Value = -2058.6667
$ws.Range("N67").Value = -5192.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3899.1667
$ws.Range("I69").Value = 3446.2856
$ws.Range("K69").Value = 10338.8568
$ws.Range("M69").Value = -9464.856800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3899.1667
$ws.Range("I72").Value = 3446.2856
$ws.Range("K72").Value = 31016.5704
$ws.Range("M72").Value = -26648.5704

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 169500.5
$ws.Range("I76").Value = 202400.6
$ws.Range("K76").Value = 202400.6
$ws.Range("M76").Value = -202085.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 169500.5
$ws.Range("I79").Value = 202400.6
$ws.Range("K79").Value = 202400.6
$ws.Range("M79").Value = -201308.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1380.68
$ws.Range("I137").Value = 956.75
$ws.Range("J137").Value = 2134.3333
$ws.Range("K137").Value = 2870.25
$ws.Range("L137").Value = 6402.999899999999
$ws.Range("M137").Value = -320.25
$ws.Range("N137").Value = -11502.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2928.611
$ws.Range("I138").Value = 6596
$ws.Range("J138").Value = 2337.0967
$ws.Range("K138").Value = 19788
$ws.Range("L138").Value = 7011.2901
$ws.Range("M138").Value = -14648
$ws.Range("N138").Value = -17291.2901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3809.11
$ws.Range("I32").Value = 2864.1604
$ws.Range("J32").Value = 7837.579
$ws.Range("K32").Value = 2864.1604
$ws.Range("L32").Value = 7837.579
$ws.Range("M32").Value = -2577.1604
$ws.Range("N32").Value = -8411.579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1385.807
$ws.Range("I61").Value = 1036.4849
$ws.Range("J61").Value = 1866.125
$ws.Range("K61").Value = 1036.4849
$ws.Range("L61").Value = 1866.125
$ws.Range("M61").Value = -824.4848999999999
$ws.Range("N61").Value = -2290.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 743.32434
$ws.Range("I74").Value = 667.37933
$ws.Range("J74").Value = 1018.625
$ws.Range("K74").Value = 667.37933
$ws.Range("L74").Value = 1018.625
$ws.Range("M74").Value = 206.62067
$ws.Range("N74").Value = -2766.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 743.32434
$ws.Range("I77").Value = 667.37933
$ws.Range("J77").Value = 1018.625
$ws.Range("K77").Value = 3336.89665
$ws.Range("L77").Value = 5093.125
$ws.Range("M77").Value = 1031.10335
$ws.Range("N77").Value = -13829.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1385.807
$ws.Range("I136").Value = 1036.4849
$ws.Range("J136").Value = 1866.125
$ws.Range("K136").Value = 3109.4547
$ws.Range("L136").Value = 5598.375
$ws.Range("M136").Value = -559.4546999999998
$ws.Range("N136").Value = -10698.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 346.42856
$ws.Range("I22").Value = 346.42856
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 346.42856
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -173.42856
$ws.Range("N22").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 65500
$ws.Range("J132").Value = 65500
$ws.Range("L132").Value = 65500
$ws.Range("N132").Value = -75620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2690136
$ws.Range("I31").Value = 1473.525
$ws.Range("J31").Value = 7578613.5
$ws.Range("K31").Value = 1473.525
$ws.Range("L31").Value = 7578613.5
$ws.Range("M31").Value = -1178.525
$ws.Range("N31").Value = -7579203.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2690136
$ws.Range("I34").Value = 1473.525
$ws.Range("J34").Value = 7578613.5
$ws.Range("K34").Value = 1473.525
$ws.Range("L34").Value = 7578613.5
$ws.Range("M34").Value = -1271.525
$ws.Range("N34").Value = -7579017.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1407.9032
$ws.Range("I58").Value = 1001.875
$ws.Range("K58").Value = 1001.875
$ws.Range("M58").Value = -798.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 41667056
$ws.Range("J94").Value = 440.41177
$ws.Range("L94").Value = 440.41177
$ws.Range("N94").Value = -1342.41177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H123").Value = 23260
$ws.Range("I123").Value = 9000
$ws.Range("K123").Value = 9000
$ws.Range("M123").Value = -4100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 970.9459000000001
$ws.Range("I134").Value = 869.7778
$ws.Range("J134").Value = 1244.1
$ws.Range("K134").Value = 2609.3334
$ws.Range("L134").Value = 3732.3
$ws.Range("M134").Value = -74.33339999999998
$ws.Range("N134").Value = -8802.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1407.9032
$ws.Range("I136").Value = 1001.875
$ws.Range("K136").Value = 3005.625
$ws.Range("M136").Value = -455.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 9960
$ws.Range("J33").Value = 87.8
$ws.Range("L33").Value = 526.8
$ws.Range("N33").Value = -1092.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 8338
$ws.Range("I63").Value = 5076
$ws.Range("J63").Value = 11600
$ws.Range("K63").Value = 15228
$ws.Range("L63").Value = 34800
$ws.Range("M63").Value = -14479
$ws.Range("N63").Value = -36298

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 8338
$ws.Range("I66").Value = 5076
$ws.Range("J66").Value = 11600
$ws.Range("K66").Value = 45684
$ws.Range("L66").Value = 104400
$ws.Range("M66").Value = -41940
$ws.Range("N66").Value = -111888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6350
$ws.Range("I80").Value = 4600
$ws.Range("J80").Value = 8100
$ws.Range("K80").Value = 4600
$ws.Range("L80").Value = 8100
$ws.Range("M80").Value = -3602
$ws.Range("N80").Value = -10096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6350
$ws.Range("I83").Value = 4600
$ws.Range("J83").Value = 8100
$ws.Range("K83").Value = 23000
$ws.Range("L83").Value = 40500
$ws.Range("M83").Value = -18008
$ws.Range("N83").Value = -50484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 47622770
$ws.Range("I122").Value = 166672830
$ws.Range("J122").Value = 2740
$ws.Range("K122").Value = 500018490
$ws.Range("L122").Value = 8220
$ws.Range("M122").Value = -500016040
$ws.Range("N122").Value = -13120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 34438.13
$ws.Range("I132").Value = 51874.2
$ws.Range("J132").Value = 2736.182
$ws.Range("K132").Value = 155622.6
$ws.Range("L132").Value = 8208.545999999998
$ws.Range("M132").Value = -153092.6
$ws.Range("N132").Value = -13268.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1029.7
$ws.Range("I82").Value = 1218.5
$ws.Range("J82").Value = 274.5
$ws.Range("K82").Value = 1218.5
$ws.Range("L82").Value = 274.5
$ws.Range("M82").Value = -857.5
$ws.Range("N82").Value = -996.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1029.7
$ws.Range("I85").Value = 1218.5
$ws.Range("J85").Value = 274.5
$ws.Range("K85").Value = 1218.5
$ws.Range("L85").Value = 274.5
$ws.Range("M85").Value = 29.5
$ws.Range("N85").Value = -2770.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2766.6667
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3150
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3150
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 2766.6667
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3150
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 15750
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8384.200000000001
$ws.Range("I126").Value = 9855.166999999999
$ws.Range("K126").Value = 29565.501
$ws.Range("M126").Value = -27095.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1133.3556
$ws.Range("I132").Value = 1001.40845
$ws.Range("J132").Value = 1626.421
$ws.Range("K132").Value = 3004.22535
$ws.Range("L132").Value = 4879.263
$ws.Range("M132").Value = -474.2253500000002
$ws.Range("N132").Value = -9939.262999999999
